$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 held a mis-tagged duplicate (GB021); re-tag it as GB020 and flip its value to 1
$ws.Range("A3").Value = "GB020"
$ws.Range("C3").Value = 1

# Rows 5 and 6 (GB023 / GB172) duplicated sources already captured above; remove them,
# shifting the remaining (empty) rows up.
$ws.Range("A5:E6").EntireRow.Delete()

# Leave the selection where the edit happened
$ws.Range("C3").Select() | Out-Null
